# Fruta / hortaliza, semanal
# Append two new data rows (62 and 63) to Sheet1, mirroring the existing
# row layout/columns A:T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the same date number format as the other "Fecha" cells (column D).
$dateFormat = $ws.Range("D61").NumberFormat

$rows = @(
    @{
        Row = 62
        A = 10
        B = "Vega Modelo de Temuco"
        C = "La Araucanía"
        D = 44890
        E = 9
        F = "Fruta"
        G = 100103
        H = "Frutos de hueso (carozo)"
        I = 100103003
        J = "Damasco"
        K = "Castle Brite"
        L = "Primera"
        M = 80
        N = 20000
        O = 20000
        P = 20000
        Q = "`$/bandeja 10 kilos"
        R = "Provincia de Limarí"
        S = 2000
        T = 10
    },
    @{
        Row = 63
        A = 10
        B = "Vega Modelo de Temuco"
        C = "La Araucanía"
        D = 44890
        E = 9
        F = "Fruta"
        G = 100103
        H = "Frutos de hueso (carozo)"
        I = 100103003
        J = "Damasco"
        K = "Castle Brite"
        L = "Primera"
        M = 80
        N = 22000
        O = 22000
        P = 22000
        Q = "`$/caja 15 kilos"
        R = "Provincia de Limarí"
        S = 1467
        T = 15
    }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
    $ws.Cells.Item($rowNum, 19).Value = $r.S
    $ws.Cells.Item($rowNum, 20).Value = $r.T
}
